# Applies updated transition probabilities to the "Park (AZ)_B" team-specific
# matrix sheet: additional simulated games changed several state-transition
# frequencies from 0 to their newly observed fractions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2857142857142857
$ws.Range("C2").Value = 0.2857142857142857
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.2857142857142857

$ws.Range("P3").Value = 1

$ws.Range("F6").Value = 0.1111111111111111
$ws.Range("J6").Value = 0.2222222222222222
$ws.Range("Q6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.5555555555555556

$ws.Range("F7").Value = 0.2
$ws.Range("Q7").Value = 0.6
$ws.Range("S7").Value = 0.2

$ws.Range("B8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.05263157894736842
$ws.Range("Q8").Value = 0.1578947368421053
$ws.Range("R8").Value = 0.05263157894736842
$ws.Range("S8").Value = 0.6842105263157895

$ws.Range("F9").Value = 0.07692307692307693
$ws.Range("Q9").Value = 0.07692307692307693
$ws.Range("R9").Value = 0.07692307692307693
$ws.Range("S9").Value = 0.7692307692307693

$ws.Range("B10").Value = 0.1
$ws.Range("F10").Value = 0.075
$ws.Range("J10").Value = 0.1
$ws.Range("O10").Value = 0.05
$ws.Range("Q10").Value = 0.075
$ws.Range("R10").Value = 0.05
$ws.Range("S10").Value = 0.55

$ws.Range("L11").Value = 1

$ws.Range("G12").Value = 1

$ws.Range("G13").Value = 1

$ws.Range("F15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.2857142857142857
$ws.Range("J15").Value = 0.2857142857142857
$ws.Range("S15").Value = 0.2857142857142857

$ws.Range("J16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.3333333333333333

$ws.Range("H17").Value = 0.3
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.5
$ws.Range("O17").Value = 0.1

$ws.Range("H18").Value = 0.25
$ws.Range("J18").Value = 0.25
$ws.Range("S18").Value = 0.5

$ws.Range("F19").Value = 0.01694915254237288
$ws.Range("H19").Value = 0.2542372881355932
$ws.Range("I19").Value = 0.1694915254237288
$ws.Range("J19").Value = 0.4067796610169492
$ws.Range("K19").Value = 0.03389830508474576
$ws.Range("M19").Value = 0.03389830508474576
$ws.Range("O19").Value = 0.05084745762711865
$ws.Range("S19").Value = 0.03389830508474576
